$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36/37: coin order swapped (TrustWalletToken <-> ImmutableX) ---
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.640'
$ws.Range("E36").Value = '  -4.97%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").Value = '  -0.81%  '

# --- Row 43/44: coin order swapped (MXToken <-> ARBITRUM) ---
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '0.945'
$ws.Range("E43").Value = '  -1.84%  '

$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D44").Value = '2.79'
$ws.Range("E44").Value = '  -0.97%  '

# --- Price / volume refresh for the remaining rows ---
$ws.Range("D2").Value = '34.488.22'
$ws.Range("E2").Value = '  +0.41%  '
$ws.Range("D3").Value = '1.808.73'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '225.36'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").Value = '0.588'
$ws.Range("E6").Value = '  +2.31%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '38.14'
$ws.Range("E8").Value = '  +5.64%  '
$ws.Range("E9").Value = '  -4.35%  '
$ws.Range("D10").Value = '0.0673'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("D11").Value = '0.0974'
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").Value = '2.069.71'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '11.12'
$ws.Range("E13").Value = '  -5.39%  '
$ws.Range("D14").Value = '1.832.24'
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("D15").Value = '34.453.53'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("E16").Value = '  -2.47%  '
$ws.Range("D17").Value = '4.41'
$ws.Range("E17").Value = '  -1.97%  '
$ws.Range("D18").Value = '68.01'
$ws.Range("D19").Value = '242.22'
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("E20").Value = '  -3.11%  '
$ws.Range("D21").Value = '11.12'
$ws.Range("E21").Value = '  -4.30%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.10'
$ws.Range("E23").Value = '  -1.87%  '
$ws.Range("E24").Value = '  +3.04%  '
$ws.Range("D25").Value = '169.89'
$ws.Range("E25").Value = '  -1.14%  '
$ws.Range("D26").Value = '7.72'
$ws.Range("E26").Value = '  -3.15%  '
$ws.Range("D27").Value = '17.55'
$ws.Range("E27").Value = '  +3.93%  '
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -0.83%  '
$ws.Range("D31").Value = '3.77'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("E33").Value = '  -4.85%  '
$ws.Range("E34").Value = '  -0.54%  '
$ws.Range("D35").Value = '1.343.07'
$ws.Range("E35").Value = '  -3.41%  '
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").Value = '  -5.60%  '
$ws.Range("E40").Value = '  +1.44%  '
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("D42").Value = '81.65'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D45").Value = '13.66'
$ws.Range("E45").Value = '  +1.58%  '
$ws.Range("E46").Value = '  +1.73%  '
$ws.Range("D47").Value = '1.970.30'
$ws.Range("E47").Value = '  +0.23%  '
$ws.Range("D48").Value = '5.74'
$ws.Range("E48").Value = '  -4.57%  '
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").Value = '102.16'
$ws.Range("E50").Value = '  -2.31%  '
$ws.Range("D51").Value = '0.0₆0120'
$ws.Range("E51").Value = '  -5.81%  '
